# Kafka Script has been added
# Update the TRANSACTION ID (column A) and TRANSACTION TIME (column D) values
# for rows 2-6 on the "Transaction Details" sheet with new, unique values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaction Details")

$ws.Range("A2").Value = "TR20240730023002227"
$ws.Range("D2").Value = "02:30:02:228"

$ws.Range("A3").Value = "TR20240730023002340"
$ws.Range("D3").Value = "02:30:02:340"

$ws.Range("A4").Value = "TR20240730023002445"
$ws.Range("D4").Value = "02:30:02:446"

$ws.Range("A5").Value = "TR20240731023002557"
$ws.Range("D5").Value = "02:30:02:557"

$ws.Range("A6").Value = "TR20240731023002680"
$ws.Range("D6").Value = "02:30:02:680"
